# Book of Gold Multichance: add a "Meta description" paragraph right
# after the title, and turn the trailing duplicate title/description
# paragraphs into a single image-generation-prompt paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new, plain (Normal-style) empty paragraph right after the
#    Heading1 title paragraph. This will become the "Meta description"
#    paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# ---------------------------------------------------------------------
# 2. The document already ends with a duplicate bold paragraph reading
#    "Play Book of Gold Multichance for Free: Review and Features"
#    (same run layout we want: empty run + bold run). Cut it from the
#    end and paste it into the new empty paragraph so we inherit its
#    exact run/formatting structure instead of re-building it by hand.
# ---------------------------------------------------------------------
$dupTitlePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$dupTitlePara.Range.Cut() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Paste() | Out-Null

# ---------------------------------------------------------------------
# 3. Re-point the pasted bold run's text from the old title to
#    "Meta description" (keep the bold formatting), then append a
#    plain run with the description text after it.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$boldTextRange = $d.Range($metaRange.Start, $metaRange.End - 1)
$boldTextRange.Text = "Meta description"

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$afterBold = $d.Range($metaRange.End - 1, $metaRange.End - 1)
$afterBold.InsertAfter(": Explore ancient Egypt with Book of Gold Multichance. Features include high volatility, free spins, and beautiful design. Play for free now.") | Out-Null

# ---------------------------------------------------------------------
# 4. The trailing italic paragraph (formerly the meta-description text,
#    now the last paragraph in the document after the cut above) gets
#    its text swapped for the new image-generation prompt, keeping its
#    italic formatting and leading empty run untouched.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastTextRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTextRange.Text = "Create a feature image for Book of Gold Multichance that showcases the game's Egyptian theme and features a happy Maya warrior with glasses. The image should be in cartoon style and use bold, bright colors to draw attention to the slot game. Incorporate symbols from the game such as the pharaoh, eye of Ra, and Ankh to give players an idea of what to expect. The Maya warrior should be depicted with a big smile and wearing glasses to symbolize the intelligence and luck required to win big in this game. Consider adding some hieroglyphics or a pyramid in the background to add to the Egyptian theme. The image should be eye-catching and memorable to help draw in players to try out Book of Gold Multichance."

Write-Output "Done"
